# "checking relay - problems with initial condition open"
# - Blank out the data rows for relays #2 and #3 (rows 3 and 4), leaving just
#   the slip-voltage/angle/frequency/delay columns (A, C:F) as empty templates
#   like the already-blank rows below them.
# - Remove the two spare template rows (8 and 9) that trailed the table.
# - Flip CB Initial State (H2) from closed (1) to open (0).
# - Move the selection to H3 and nudge the license textbox up to sit right
#   under the (now shorter) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CB Initial State for relay #1: closed -> open
$ws.Range("H2").Value = 0

# Relay #2 (row 3) and #3 (row 4): clear out all the per-relay settings,
# leaving only the empty slip-voltage/angle/frequency/delay cells (matching
# the blank template rows 5-7).
foreach ($r in 3, 4) {
    $ws.Range("B$r").Clear()
    $ws.Range("G$r`:V$r").Clear()
    $ws.Range("A$r").ClearContents()
    $ws.Range("C$r`:F$r").ClearContents()
}

# Drop the two trailing spare rows (8 and 9).
$ws.Rows("8:9").Delete()

# Move the license/comment textbox up to row 5 (0-based row 4) now that the
# table is shorter. Keep the same vertical offset within the row (164523 EMU
# = 12.954... pt) that it already had relative to its anchor row.
$shp = $ws.Shapes.Item("TextBox 3")
$rowOffPt = 164523 / 12700
$shp.Top = $ws.Rows.Item(5).Top + $rowOffPt

# Update the active selection.
$null = $ws.Range("H3").Select()
